# Update workbook for the 09-01-2025 daily report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header date labels -------------------------------------------------
$ws.Range("A1").Value = "Mangrove Communication   09.01.2025"
$ws.Range("A10").Value = "DAILY STOCK                         (09/01/2025) "

# --- Sale/stock table (rows 3-6) ----------------------------------------
$ws.Range("C3").Value = 29275
$ws.Range("D3").Value = 3199

$ws.Range("C4").Value = 45936
$ws.Range("D4").Value = 15244

$ws.Range("C5").Value = 35540
$ws.Range("D5").Value = 10257
$ws.Range("F5").Value = $null

$ws.Range("C6").Value = 64350
$ws.Range("D6").Value = 11668

# --- I top up balance table (rows 13-27) ---------------------------------
$ws.Range("C13").Value = 145970

$ws.Range("C14").Value = 240567
$ws.Range("D14").Value = 175101
$ws.Range("E14").Value = 432107

$ws.Range("C18").Value = 320

$ws.Range("D20").Value = 2100
$ws.Range("E20").Value = 5000

$ws.Range("C21").Value = 230
$ws.Range("D21").Value = 30

$ws.Range("C22").Value = 1000
$ws.Range("E22").Value = 1000

$ws.Range("C24").Value = 40

$ws.Range("C26").Value = 82
$ws.Range("D26").Value = 46

$ws.Range("C27").Value = 72

# --- View state (best effort) --------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("K28").Select()
